$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "survey" sheet: insert a new row 21 for the "maintenance_priority" field
#    (select_one service_priority_list) right before the existing
#    "voltage_regulator" row, which shifts down to row 22 (and the trailing
#    "end screen" row shifts from 22 to 23).
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Rows.Item(21).Insert()
$survey.Rows.Item(21).RowHeight = 15

$survey.Cells.Item(21, 4).Value  = "select_one"
$survey.Cells.Item(21, 5).Value  = "service_priority_list"
$survey.Cells.Item(21, 6).Value  = "maintenance_priority"
$survey.Cells.Item(21, 7).Value  = "Service Requested With Priority"
$survey.Cells.Item(21, 8).Value  = "Servicio solicitado con prioridad"
$survey.Cells.Item(21, 9).Value  = "Choose priority for maintenance if applicable:"
$survey.Cells.Item(21, 10).Value = "Seleccione la prioridad de mantenimiento si corresponde:"

# Widen column J (name column) slightly to fit the new content.
$survey.Columns.Item(10).ColumnWidth = 34

# ---------------------------------------------------------------------------
# 2. "choices" sheet: add the new "service_priority_list" choice list
#    (low / medium / high / not_applicable) as rows 59-62.
#    Values are written column-by-column to match shared-string ordering.
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Cells.Item(59, 1).Value = "service_priority_list"
$choices.Cells.Item(60, 1).Value = "service_priority_list"
$choices.Cells.Item(61, 1).Value = "service_priority_list"
$choices.Cells.Item(62, 1).Value = "service_priority_list"

$choices.Cells.Item(59, 2).Value = "low"
$choices.Cells.Item(60, 2).Value = "medium"
$choices.Cells.Item(61, 2).Value = "high"
$choices.Cells.Item(62, 2).Value = "not_applicable"

$choices.Cells.Item(59, 3).Value = "Low"
$choices.Cells.Item(60, 3).Value = "Medium"
$choices.Cells.Item(61, 3).Value = "High"
$choices.Cells.Item(62, 3).Value = "Not Applicable"

$choices.Cells.Item(59, 4).Value = "Bajo"
$choices.Cells.Item(60, 4).Value = "Medio"
$choices.Cells.Item(61, 4).Value = "Alto"
$choices.Cells.Item(62, 4).Value = "No Aplica"

# ---------------------------------------------------------------------------
# 3. Restore per-sheet selections.
# ---------------------------------------------------------------------------
$survey.Range("E21").Select()

# "choices" becomes the active/selected tab, with B67 selected.
$choices.Activate()
$choices.Range("B67").Select()
